$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 112, shifting existing rows 112:203 down to 113:204
$ws.Rows.Item(112).Insert()

# Populate the new row 112 with the new observation data
$ws.Range("A112").Value = 11
$ws.Range("B112").Value = "Vega Monumental Concepción"
$ws.Range("C112").Value = "Bíobío"
$ws.Range("D112").Value = 44651
$ws.Range("D112").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E112").Value = 8
$ws.Range("F112").Value = 100114013
$ws.Range("G112").Value = "Zanahoria"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 220
$ws.Range("K112").Value = 6000
$ws.Range("L112").Value = 7000
$ws.Range("M112").Value = 6545
$ws.Range("N112").Value = "`$/saco 20 kilos"
$ws.Range("O112").Value = "Chillán"
$ws.Range("P112").Value = 327
$ws.Range("Q112").Value = 20
$ws.Range("R112").Value = "Hortaliza"
